# Updated symbol list on Wed Dec 21 18:48:34 UTC 2022 with GitHub Actions
#
# This script updates the "Price" (column D) values for most rows in the
# crypto listing sheet, and also reshuffles the order of three rows
# (BKEXToken / CEJI / KickToken) along with updating their price / link /
# text fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        $Worksheet,
        [string]$Address,
        [string]$Value
    )
    $range = $Worksheet.Range($Address)
    # Force a Text number format so that numeric-looking strings (prices)
    # are preserved exactly (including trailing/leading zeros) instead of
    # being coerced into a floating point number by Excel.
    $range.NumberFormat = "@"
    $range.Value = $Value
}

function Set-PlainCell {
    param(
        $Worksheet,
        [string]$Address,
        [string]$Value
    )
    $Worksheet.Range($Address).Value = $Value
}

# --- Column D (Price) straightforward updates -----------------------------
Set-TextCell $ws "D2"  "247.65"
Set-TextCell $ws "D3"  "22.31"
Set-TextCell $ws "D4"  "5.247"
Set-TextCell $ws "D5"  "0.05686"
Set-TextCell $ws "D7"  "6.314"
Set-TextCell $ws "D8"  "0.8073"
Set-TextCell $ws "D9"  "0.8624"
Set-TextCell $ws "D10" "0.1415"
Set-TextCell $ws "D11" "0.07392"
Set-TextCell $ws "D12" "0.03049"
Set-TextCell $ws "D13" "0.03078"
Set-TextCell $ws "D14" "0.09391"
Set-TextCell $ws "D15" "3.880"
Set-TextCell $ws "D16" "0.001581"
Set-TextCell $ws "D17" "0.04782"
Set-TextCell $ws "D18" "0.01829"

# Row 19 also gets an extra "Worstin24h" suffix on its E column text.
Set-TextCell  $ws "D19" "0.0005802"
Set-PlainCell $ws "E19" "18OneONEWorstin24h"

Set-TextCell $ws "D20" "0.006439"
Set-TextCell $ws "D21" "0.005035"
Set-TextCell $ws "D22" "0.0009967"
Set-TextCell $ws "D23" "0.0001501"
Set-TextCell $ws "D24" "3.691"
Set-TextCell $ws "D25" "2.195"
Set-TextCell $ws "D26" "0.3248"
Set-TextCell $ws "D27" "0.1351"

Set-TextCell $ws "D40" "0.03958"

# --- Rows 41-43: rows reshuffled (BKEXToken / CEJI / KickToken) -----------
# Before: 41=BKEXToken, 42=CEJI, 43=KickToken
# After:  41=KickToken, 42=BKEXToken, 43=CEJI

Set-PlainCell $ws "B41" "KickToken"
Set-PlainCell $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell  $ws "D41" "0.006815"
Set-PlainCell $ws "E41" "40KickTokenKICK"

Set-PlainCell $ws "B42" "BKEXToken"
Set-PlainCell $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell  $ws "D42" "0.1065"
Set-PlainCell $ws "E42" "41BKEXTokenBKK"

Set-PlainCell $ws "B43" "CEJI"
Set-PlainCell $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell  $ws "D43" "0.003202"
Set-PlainCell $ws "E43" "42CEJICEJI"

# --- Remaining column D updates -------------------------------------------
Set-TextCell $ws "D44" "0.008466"
Set-TextCell $ws "D45" "0.00005595"
Set-TextCell $ws "D47" "0.4501"
Set-TextCell $ws "D48" "0.1959"
Set-TextCell $ws "D50" "0.01011"
